# Remove the "BrowserProcessName" config row from the Settings sheet.
# This deletes the entire worksheet row 14 (Name=BrowserProcessName,
# Value=chrome.exe, Description=...), shifting all following rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

$ws.Rows.Item(14).Delete()
